$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")
$ws.Range("A484").Value = "honey_block"
$ws.Range("E484").Value = "1.15"
$ws.Range("F484").Value = "slime_block"
$ws.Range("A485").Value = "honeycomb_block"
$ws.Range("E485").Value = "1.15"
$ws.Range("F485").Value = "orange_terracotta"
$ws.Range("A486").Value = "ancient_debris"
$ws.Range("E486").Value = "1.16"
$ws.Range("F486").Value = "netherrack"
$ws.Range("A487").Value = "basalt"
$ws.Range("E487").Value = "1.16"
$ws.Range("F487").Value = "cobblestone"
$ws.Range("A488").Value = "netherite_block"
$ws.Range("E488").Value = "1.16"
$ws.Range("F488").Value = "obsidian"
$ws.Range("A489").Value = "crimson_stem"
$ws.Range("B489").Value = "???"
$ws.Range("E489").Value = "1.16"
$ws.Range("F489").Value = "spruce_log"
$ws.Range("A490").Value = "warped_stem"
$ws.Range("B490").Value = "???"
$ws.Range("E490").Value = "1.16"
$ws.Range("F490").Value = "spruce_log"
$ws.Range("A491").Value = "crimson_planks"
$ws.Range("E491").Value = "1.16"
$ws.Range("F491").Value = "spruce_planks"
$ws.Range("A492").Value = "warped_planks"
$ws.Range("E492").Value = "1.16"
$ws.Range("F492").Value = "spruce_planks"
$ws.Range("A493").Value = "crimson_sign"
$ws.Range("B493").Value = "???"
$ws.Range("E493").Value = "1.16"
$ws.Range("F493").Value = "spruce_sign"
$ws.Range("A494").Value = "warped_sign"
$ws.Range("B494").Value = "???"
$ws.Range("E494").Value = "1.16"
$ws.Range("F494").Value = "spruce_sign"
$ws.Range("A495").Value = "crimson_wall_sign"
$ws.Range("B495").Value = "???"
$ws.Range("E495").Value = "1.16"
$ws.Range("F495").Value = "spruce_wall_sign"
$ws.Range("A496").Value = "warped_wall_sign"
$ws.Range("B496").Value = "???"
$ws.Range("E496").Value = "1.16"
$ws.Range("F496").Value = "spruce_wall_sign"
$ws.Range("A497").Value = "crimson_slab"
$ws.Range("B497").Value = "???"
$ws.Range("E497").Value = "1.16"
$ws.Range("F497").Value = "spruce_slab"
$ws.Range("A498").Value = "warped_slab"
$ws.Range("B498").Value = "???"
$ws.Range("E498").Value = "1.16"
$ws.Range("F498").Value = "spruce_slab"
$ws.Range("A499").Value = "crimson_fence"
$ws.Range("B499").Value = "???"
$ws.Range("E499").Value = "1.16"
$ws.Range("F499").Value = "spruce_fence"
$ws.Range("A500").Value = "warped_fence"
$ws.Range("B500").Value = "???"
$ws.Range("E500").Value = "1.16"
$ws.Range("F500").Value = "spruce_fence"
$ws.Range("A501").Value = "crimson_fence_gate"
$ws.Range("B501").Value = "???"
$ws.Range("E501").Value = "1.16"
$ws.Range("F501").Value = "spruce_fence_gate"
$ws.Range("A502").Value = "warped_fence_gate"
$ws.Range("B502").Value = "???"
$ws.Range("E502").Value = "1.16"
$ws.Range("F502").Value = "spruce_fence_gate"
$ws.Range("A503").Value = "crimson_pressure_plate"
$ws.Range("E503").Value = "1.16"
$ws.Range("F503").Value = "spruce_pressure_plate"
$ws.Range("A504").Value = "warped_pressure_plate"
$ws.Range("E504").Value = "1.16"
$ws.Range("F504").Value = "spruce_pressure_plate"
$ws.Range("A505").Value = "crimson_button"
$ws.Range("B505").Value = "???"
$ws.Range("E505").Value = "1.16"
$ws.Range("F505").Value = "spruce_button"
$ws.Range("A506").Value = "warped_button"
$ws.Range("B506").Value = "???"
$ws.Range("E506").Value = "1.16"
$ws.Range("F506").Value = "spruce_button"
$ws.Range("A507").Value = "crimson_door"
$ws.Range("B507").Value = "???"
$ws.Range("E507").Value = "1.16"
$ws.Range("F507").Value = "spruce_door"
$ws.Range("A508").Value = "warped_door"
$ws.Range("B508").Value = "???"
$ws.Range("E508").Value = "1.16"
$ws.Range("F508").Value = "spruce_door"
$ws.Range("A509").Value = "crimson_trapdoor"
$ws.Range("B509").Value = "???"
$ws.Range("E509").Value = "1.16"
$ws.Range("F509").Value = "spruce_trapdoor"
$ws.Range("A510").Value = "warped_trapdoor"
$ws.Range("B510").Value = "???"
$ws.Range("E510").Value = "1.16"
$ws.Range("F510").Value = "spruce_trapdoor"
$ws.Range("A511").Value = "crimson_fungus"
$ws.Range("E511").Value = "1.16"
$ws.Range("F511").Value = "red_mushroom"
$ws.Range("A512").Value = "warped_mushroom"
$ws.Range("E512").Value = "1.16"
$ws.Range("F512").Value = "brown_mushroom"
$ws.Range("A513").Value = "crimson_nylium"
$ws.Range("E513").Value = "1.16"
$ws.Range("F513").Value = "netherrack"
$ws.Range("A514").Value = "warped_nylium"
$ws.Range("E514").Value = "1.16"
$ws.Range("F514").Value = "netherrack"
$ws.Range("A515").Value = "crimson_roots"
$ws.Range("E515").Value = "1.16"
$ws.Range("A516").Value = "warped_roots"
$ws.Range("E516").Value = "1.16"
$ws.Range("A517").Value = "nether_sprouts"
$ws.Range("E517").Value = "1.16"
$ws.Range("A518").Value = "shroomlight"
$ws.Range("E518").Value = "1.16"
$ws.Range("F518").Value = "glowstone"
$ws.Range("A519").Value = "soul_fire"
$ws.Range("E519").Value = "1.16"
$ws.Range("F519").Value = "fire"
$ws.Range("A520").Value = "soul_lantern"
$ws.Range("B520").Value = "???"
$ws.Range("E520").Value = "1.16"
$ws.Range("F520").Value = "lantern"
$ws.Range("A521").Value = "soul_torch"
$ws.Range("E521").Value = "1.16"
$ws.Range("F521").Value = "torch"
$ws.Range("A522").Value = "soul_wall_torch"
$ws.Range("B522").Value = "???"
$ws.Range("E522").Value = "1.16"
$ws.Range("F522").Value = "wall_torch"
$ws.Range("A523").Value = "soul_soil"
$ws.Range("E523").Value = "1.16"
$ws.Range("F523").Value = "soul_sand"
$ws.Range("A524").Value = "warped_wart_block"
$ws.Range("E524").Value = "1.16"
$ws.Range("F524").Value = "nether_wart_block"
$ws.Range("A525").Value = "weeping_vines"
$ws.Range("E525").Value = "1.16"
$ws.Range("A526").Value = "crying_obsidian"
$ws.Range("E526").Value = "1.16"
$ws.Range("F526").Value = "obsidian"
$ws.Range("A527").Value = "target"
$ws.Range("E527").Value = "1.16"
$ws.Range("F527").Value = "chiseled_stone_bricks"
$ws.Range("A528").Value = "crimson_hyphae"
$ws.Range("B528").Value = "???"
$ws.Range("E528").Value = "1.16"
$ws.Range("F528").Value = "spruce_log"
$ws.Range("A529").Value = "warped_hyphae"
$ws.Range("B529").Value = "???"
$ws.Range("E529").Value = "1.16"
$ws.Range("F529").Value = "spruce_log"
$ws.Range("A530").Value = "nether_gold_ore"
$ws.Range("E530").Value = "1.16"
$ws.Range("F530").Value = "netherrack"
$ws.Range("A531").Value = "twisting_vines"
$ws.Range("E531").Value = "1.16"
$ws.Range("A532").Value = "polished_basalt"
$ws.Range("E532").Value = "1.16"
$ws.Range("F532").Value = "cobblestone"
$ws.Range("A533").Value = "respawn_anchor"
$ws.Range("E533").Value = "1.16"
$ws.Range("F533").Value = "obsidian"
$ws.Range("A534").Value = "lodestone"
$ws.Range("E534").Value = "1.16"
$ws.Range("F534").Value = "chiseled_stone_bricks"
$ws.Range("A535").Value = "blackstone"
$ws.Range("E535").Value = "1.16"
$ws.Range("F535").Value = "cobblestone"
$ws.Range("A536").Value = "blackstone_slab"
$ws.Range("E536").Value = "1.16"
$ws.Range("F536").Value = "cobblestone_slab"
$ws.Range("A537").Value = "blackstone_stairs"
$ws.Range("E537").Value = "1.16"
$ws.Range("F537").Value = "cobblestone_stairs"
$ws.Range("A538").Value = "blackstone_wall"
$ws.Range("E538").Value = "1.16"
$ws.Range("F538").Value = "cobblestone_wall"
$ws.Range("A539").Value = "polished_blackstone"
$ws.Range("E539").Value = "1.16"
$ws.Range("F539").Value = "cobblestone"
$ws.Range("A540").Value = "polished_blackstone_slab"
$ws.Range("E540").Value = "1.16"
$ws.Range("F540").Value = "cobblestone_slab"
$ws.Range("A541").Value = "polished_blackstone_stairs"
$ws.Range("E541").Value = "1.16"
$ws.Range("F541").Value = "cobblestone_stairs"
$ws.Range("A542").Value = "polished_blackstone_wall"
$ws.Range("E542").Value = "1.16"
$ws.Range("F542").Value = "cobblestone_wall"
$ws.Range("A543").Value = "polished_blackstone_bricks"
$ws.Range("E543").Value = "1.16"
$ws.Range("F543").Value = "stone_bricks"
$ws.Range("A544").Value = "polished_blackstone_brick_slab"
$ws.Range("E544").Value = "1.16"
$ws.Range("F544").Value = "stone_brick_slab"
$ws.Range("A545").Value = "polished_blackstone_brick_stairs"
$ws.Range("E545").Value = "1.16"
$ws.Range("F545").Value = "stone_brick_stairs"
$ws.Range("A546").Value = "polished_blackstone_brick_wall"
$ws.Range("E546").Value = "1.16"
$ws.Range("F546").Value = "stone_brick_wall"
$ws.Range("A547").Value = "polished_blackstone_button"
$ws.Range("B547").Value = "???"
$ws.Range("E547").Value = "1.16"
$ws.Range("F547").Value = "stone_button"
$ws.Range("A548").Value = "polished_blackstone_pressure_plate"
$ws.Range("E548").Value = "1.16"
$ws.Range("F548").Value = "stone_pressure_plate"
$ws.Range("A549").Value = "chiseled_nether_bricks"
$ws.Range("E549").Value = "1.16"
$ws.Range("F549").Value = "nether_bricks"
$ws.Range("A550").Value = "chiseled_polished_blackstone"
$ws.Range("E550").Value = "1.16"
$ws.Range("F550").Value = "chiseled_stone_bricks"
$ws.Range("A551").Value = "cracked_nether_bricks"
$ws.Range("E551").Value = "1.16"
$ws.Range("F551").Value = "nether_bricks"
$ws.Range("A552").Value = "cracked_polished_blackstone_bricks"
$ws.Range("E552").Value = "1.16"
$ws.Range("F552").Value = "cracked_stone_bricks"
$ws.Range("A553").Value = "gilded_blackstone"
$ws.Range("E553").Value = "1.16"
$ws.Range("F553").Value = "cobblestone"
$ws.Range("A554").Value = "quartz_bricks"
$ws.Range("E554").Value = "1.16"
$ws.Range("F554").Value = "nether_quartz_block"
$ws.Range("A555").Value = "soul_campfire"
$ws.Range("E555").Value = "1.16"
$ws.Range("F555").Value = "campfire"
$ws.Range("A556").Value = "chains"
$ws.Range("E556").Value = "1.16"
$ws.Range("F556").Value = "iron_bars"
$ws.Range("E557").Value = "1.17"
$ws.Range("E558").Value = "1.17"
$ws.Range("E559").Value = "1.17"
$ws.Range("E560").Value = "1.17"
$ws.Range("E561").Value = "1.17"
$ws.Range("E562").Value = "1.17"
$ws.Range("E563").Value = "1.17"
$ws.Range("E564").Value = "1.17"
$ws.Range("E565").Value = "1.17"
$ws.Range("E566").Value = "1.17"
$ws.Range("E567").Value = "1.17"
$ws.Range("E568").Value = "1.17"
$ws.Range("E569").Value = "1.17"
$ws.Range("E570").Value = "1.17"
$ws.Range("E571").Value = "1.17"
$ws.Range("E572").Value = "1.17"
$ws.Range("E573").Value = "1.17"
$ws.Range("E574").Value = "1.17"
$ws.Range("E575").Value = "1.17"
$ws.Range("E576").Value = "1.17"
$ws.Range("E577").Value = "1.17"
$ws.Range("E578").Value = "1.17"
$ws.Range("E579").Value = "1.17"
$ws.Range("E580").Value = "1.17"
$ws.Range("E581").Value = "1.17"
$ws.Range("E582").Value = "1.17"
$ws.Range("E583").Value = "1.17"
$ws.Range("E584").Value = "1.17"
$ws.Range("E585").Value = "1.17"
$ws.Range("E586").Value = "1.17"
$ws.Range("E587").Value = "1.17"
$ws.Range("E588").Value = "1.17"
$ws.Range("E589").Value = "1.17"
$ws.Range("E590").Value = "1.17"
$ws.Range("E591").Value = "1.17"
$ws.Range("E592").Value = "1.17"
$ws.Range("E593").Value = "1.17"
$ws.Range("E594").Value = "1.17"
$ws.Range("E595").Value = "1.17"
$ws.Range("E596").Value = "1.17"
$ws.Range("E597").Value = "1.17"
$ws.Range("E598").Value = "1.17"
$ws.Range("E599").Value = "1.17"
$ws.Range("E600").Value = "1.17"
$ws.Range("E601").Value = "1.17"
$ws.Range("E602").Value = "1.17"
$ws.Range("E603").Value = "1.17"
$ws.Range("E604").Value = "1.17"
$ws.Range("E605").Value = "1.17"
$ws.Range("E606").Value = "1.17"
$ws.Range("E607").Value = "1.17"
$ws.Range("E608").Value = "1.17"
$ws.Activate()
$ws.Range("F557").Select()
